# Update countries & provincias Spain
# Refresh COVID-19 stats table ("Pais" sheet) with the latest pull:
#  - updates case counts for several countries (re-sorted by "Casos totales")
#  - several rows now show a different country because the table is sorted
#    by total cases descending and rankings shifted with the new data
#  - refreshes the "Datos actualizados ..." timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 23 de Septiembre de 2020 a las 01:52"

# Row 4
$ws.Cells.Item(4, 2).Value = 7096438
$ws.Cells.Item(4, 3).Value = 34197
$ws.Cells.Item(4, 4).Value = 4343501
$ws.Cells.Item(4, 5).Value = 2547533
$ws.Cells.Item(4, 7).Value = 902
$ws.Cells.Item(4, 8).Value = 205404

# Row 6
$ws.Cells.Item(6, 2).Value = 4595335
$ws.Cells.Item(6, 3).Value = 35252
$ws.Cells.Item(6, 5).Value = 511549
$ws.Cells.Item(6, 7).Value = 809
$ws.Cells.Item(6, 8).Value = 138159

# Row 25
$ws.Cells.Item(25, 4).Value = 247900
$ws.Cells.Item(25, 5).Value = 19785

# Row 39
$ws.Cells.Item(39, 2).Value = 102254
$ws.Cells.Item(39, 3).Value = 113
$ws.Cells.Item(39, 4).Value = 91143
$ws.Cells.Item(39, 5).Value = 5305
$ws.Cells.Item(39, 7).Value = 19
$ws.Cells.Item(39, 8).Value = 5806

# Row 59 -> Chequia
$ws.Cells.Item(59, 1).Value = "Chequia"
$ws.Cells.Item(59, 2).Value = 53158
$ws.Cells.Item(59, 3).Value = 2394
$ws.Cells.Item(59, 4).Value = 26151
$ws.Cells.Item(59, 5).Value = 26476
$ws.Cells.Item(59, 7).Value = 9
$ws.Cells.Item(59, 8).Value = 531

# Row 60 -> Uzbekistan
$ws.Cells.Item(60, 1).Value = "Uzbekistan"
$ws.Cells.Item(60, 2).Value = 52685
$ws.Cells.Item(60, 3).Value = 615
$ws.Cells.Item(60, 4).Value = 49067
$ws.Cells.Item(60, 5).Value = 3176
$ws.Cells.Item(60, 7).Value = 5
$ws.Cells.Item(60, 8).Value = 442

# Row 121 -> Nicaragua
$ws.Cells.Item(121, 1).Value = "Nicaragua"
$ws.Cells.Item(121, 2).Value = 5073
$ws.Cells.Item(121, 3).Value = 112
$ws.Cells.Item(121, 4).Value = 2913
$ws.Cells.Item(121, 5).Value = 2011
$ws.Cells.Item(121, 7).Value = 2
$ws.Cells.Item(121, 8).Value = 149

# Row 122 -> Hong Kong
$ws.Cells.Item(122, 1).Value = "Hong Kong"
$ws.Cells.Item(122, 2).Value = 5047
$ws.Cells.Item(122, 3).Value = 8
$ws.Cells.Item(122, 4).Value = 4717
$ws.Cells.Item(122, 5).Value = 227
$ws.Cells.Item(122, 8).Value = 103

# Row 123 -> Guinea Ecuatorial
$ws.Cells.Item(123, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(123, 2).Value = 5018
$ws.Cells.Item(123, 3).Value = 16
$ws.Cells.Item(123, 4).Value = 4509
$ws.Cells.Item(123, 5).Value = 426
$ws.Cells.Item(123, 8).Value = 83

# Row 124 -> Congo
$ws.Cells.Item(124, 1).Value = "Congo"
$ws.Cells.Item(124, 2).Value = 5002
$ws.Cells.Item(124, 4).Value = 3887
$ws.Cells.Item(124, 5).Value = 1026
$ws.Cells.Item(124, 8).Value = 89

# Row 126
$ws.Cells.Item(126, 2).Value = 4759
$ws.Cells.Item(126, 3).Value = 19
$ws.Cells.Item(126, 4).Value = 4541
$ws.Cells.Item(126, 7).Value = 3
$ws.Cells.Item(126, 8).Value = 100

# Row 129 -> Guadalupe
$ws.Cells.Item(129, 1).Value = "Guadalupe"
$ws.Cells.Item(129, 2).Value = 4487
$ws.Cells.Item(129, 3).Value = 1061
$ws.Cells.Item(129, 4).Value = 2199
$ws.Cells.Item(129, 5).Value = 2246
$ws.Cells.Item(129, 7).Value = 16
$ws.Cells.Item(129, 8).Value = 42

# Row 130 -> Angola
$ws.Cells.Item(130, 1).Value = "Angola"
$ws.Cells.Item(130, 2).Value = 4236
$ws.Cells.Item(130, 3).Value = 119
$ws.Cells.Item(130, 4).Value = 1462
$ws.Cells.Item(130, 5).Value = 2619
$ws.Cells.Item(130, 7).Value = 1
$ws.Cells.Item(130, 8).Value = 155

# Row 131 -> Trinidad yTobago
$ws.Cells.Item(131, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(131, 2).Value = 4026
$ws.Cells.Item(131, 3).Value = 81
$ws.Cells.Item(131, 4).Value = 1871
$ws.Cells.Item(131, 5).Value = 2090
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 65

# Row 132 -> Georgia
$ws.Cells.Item(132, 1).Value = "Georgia"
$ws.Cells.Item(132, 2).Value = 3913
$ws.Cells.Item(132, 3).Value = 218
$ws.Cells.Item(132, 4).Value = 1574
$ws.Cells.Item(132, 5).Value = 2316
$ws.Cells.Item(132, 8).Value = 23

# Row 133 -> Siria
$ws.Cells.Item(133, 1).Value = "Siria"
$ws.Cells.Item(133, 2).Value = 3877
$ws.Cells.Item(133, 3).Value = 44
$ws.Cells.Item(133, 4).Value = 983
$ws.Cells.Item(133, 5).Value = 2716
$ws.Cells.Item(133, 7).Value = 3
$ws.Cells.Item(133, 8).Value = 178

# Row 134 -> Lituania
$ws.Cells.Item(134, 1).Value = "Lituania"
$ws.Cells.Item(134, 2).Value = 3859
$ws.Cells.Item(134, 3).Value = 45
$ws.Cells.Item(134, 4).Value = 2225
$ws.Cells.Item(134, 5).Value = 1547
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 87

# Row 135 -> Aruba
$ws.Cells.Item(135, 1).Value = "Aruba"
$ws.Cells.Item(135, 2).Value = 3665
$ws.Cells.Item(135, 3).Value = 78
$ws.Cells.Item(135, 4).Value = 2426
$ws.Cells.Item(135, 5).Value = 1214
$ws.Cells.Item(135, 7).Value = 1
$ws.Cells.Item(135, 8).Value = 25

# Row 136 -> Mayotte
$ws.Cells.Item(136, 1).Value = "Mayotte"
$ws.Cells.Item(136, 2).Value = 3541
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 2964
$ws.Cells.Item(136, 5).Value = 537
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 40

# Row 137 -> Gambia
$ws.Cells.Item(137, 1).Value = "Gambia"
$ws.Cells.Item(137, 2).Value = 3540
$ws.Cells.Item(137, 3).Value = 14
$ws.Cells.Item(137, 4).Value = 2002
$ws.Cells.Item(137, 5).Value = 1428
$ws.Cells.Item(137, 7).Value = 2
$ws.Cells.Item(137, 8).Value = 110

# Row 138 -> Tailandia
$ws.Cells.Item(138, 1).Value = "Tailandia"
$ws.Cells.Item(138, 2).Value = 3511
$ws.Cells.Item(138, 3).Value = 5
$ws.Cells.Item(138, 4).Value = 3343
$ws.Cells.Item(138, 5).Value = 109
$ws.Cells.Item(138, 8).Value = 59

# Row 139 -> Somalia
$ws.Cells.Item(139, 1).Value = "Somalia"
$ws.Cells.Item(139, 2).Value = 3465
$ws.Cells.Item(139, 4).Value = 2877
$ws.Cells.Item(139, 5).Value = 490
$ws.Cells.Item(139, 8).Value = 98

# Row 150
$ws.Cells.Item(150, 2).Value = 2324
$ws.Cells.Item(150, 3).Value = 21
$ws.Cells.Item(150, 4).Value = 1549
$ws.Cells.Item(150, 5).Value = 736

# Row 154
$ws.Cells.Item(154, 2).Value = 1934
$ws.Cells.Item(154, 3).Value = 7
$ws.Cells.Item(154, 4).Value = 1645
$ws.Cells.Item(154, 5).Value = 243

# Row 165 -> Martinica
$ws.Cells.Item(165, 1).Value = "Martinica"
$ws.Cells.Item(165, 2).Value = 1290
$ws.Cells.Item(165, 3).Value = 168
$ws.Cells.Item(165, 4).Value = 98
$ws.Cells.Item(165, 5).Value = 1172
$ws.Cells.Item(165, 7).Value = 2
$ws.Cells.Item(165, 8).Value = 20

# Row 166 -> Niger
$ws.Cells.Item(166, 1).Value = "Niger"
$ws.Cells.Item(166, 2).Value = 1193
$ws.Cells.Item(166, 3).Value = 4
$ws.Cells.Item(166, 4).Value = 1104
$ws.Cells.Item(166, 5).Value = 20
$ws.Cells.Item(166, 8).Value = 69

# Row 167 -> Republica del Chad
$ws.Cells.Item(167, 1).Value = "Republica del Chad"
$ws.Cells.Item(167, 2).Value = 1155
$ws.Cells.Item(167, 3).Value = 2
$ws.Cells.Item(167, 4).Value = 967
$ws.Cells.Item(167, 5).Value = 107
$ws.Cells.Item(167, 8).Value = 81

# Row 174
$ws.Cells.Item(174, 2).Value = 527
$ws.Cells.Item(174, 3).Value = 10
$ws.Cells.Item(174, 5).Value = 288

# Row 214 -> Montserrat
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1

# Row 215 -> Islas Malvinas
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
